$d = $word.ActiveDocument

# The "_GoBack" bookmark currently wraps nothing at the end of the last
# paragraph ("Learning about props"). In the edited document it ends up
# after the new "Learned about implicit returns" run instead, so remove
# it from its current location; it will be re-created in the right spot
# as part of the new content inserted below.
$bookmarks = $d.Bookmarks
if ($bookmarks.Exists("_GoBack")) {
    $bookmarks.Item("_GoBack").Delete()
}

# Start a fresh, empty paragraph right after the current last paragraph
# ("Learning about props") so the new journal entries can be inserted
# into it without disturbing the existing text.
$lastParagraph = $d.Paragraphs.Last
$insertionPoint = $lastParagraph.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range

$newParagraphsXml =
    '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Learning about passing dynamic data</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>A component is basically just an object</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Console gives comments about possible </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>codesmells</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> or </w:t></w:r>' +
      '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>unused pieces of code</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Learning about stateless functional components</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Learned about implicit returns</w:t></w:r>' +
      '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

# A trailing empty paragraph (matching the diff's final blank w:p) is
# produced automatically because InsertXML replaces the lone paragraph
# mark that InsertParagraphAfter() just created, and the story needs a
# new closing mark after the inserted content - it inherits the same
# "en-GB" paragraph mark formatting the empty paragraph already had.
$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body>' + $newParagraphsXml + '</w:body>' +
          '</w:document>' +
        '</pkg:xmlData>' +
      '</pkg:part>' +
    '</pkg:package>'

$newRange.InsertXML($packageXml) | Out-Null
